$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the shared-string text used by column B rows 133-212:
#    "Volumnometric water content" -> "Volumetric water content " (typo fix + trailing space)
$newText = "Volumetric water content "
for ($r = 133; $r -le 212; $r++) {
    $ws.Range("B$r").Value = $newText
}

# 2) Turn on word-wrap for that same B133:B212 range (adds a second cellXfs entry)
$ws.Range("B133:B212").WrapText = $true

# 3) Update the view state: selection moves to B215 (and the view had scrolled so
#    that A195 is the top-left visible cell)
$excel.ActiveWindow.ScrollRow = 195
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B215").Select() | Out-Null
